$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row -> (new C value, new E value)
# C = nombre_aides, E = montant_total
$updates = @{
    2   = @(766329, 1429227365)
    93  = @(16941, 50660115)
    100 = @(9344, 23829850)
    115 = @(81805, 436644943)
    121 = @(1306282, 2275122346)
    129 = @(633627, 3432082404)
    130 = @(4247, 141297718)
    132 = @(585880, 3468923168)
    136 = @(26695, 144329346)
    178 = @(515885, 891200284)
    237 = @(283320, 1438419995)
    240 = @(205917, 1069507520)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}
